$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header date text in B1 (kept as plain text, not a date serial)
$ws.Range("B1").Value = "24/03/2023"

# Column C: update all rows 2-20 to 10
$ws.Range("C2:C20").Value = 10

# Column B: new values for rows 2-20
$bValues = @(270.4, 159, 2, 73, 29, 50, 29, 71, 293, 204, 331.5, 233, 354, 113, 151, 47, 1, 6, 28)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
